$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Acknowledgments")

# Find the row containing "lxml" in column A (the package name column)
$lxmlCell = $ws.Range("A1:A1000").Find("lxml", [Type]::Missing, [Type]::Missing, 1)
if ($lxmlCell -ne $null) {
    $rowNum = $lxmlCell.Row
    $ws.Rows.Item($rowNum).Delete()
}
